# Update TIP contact status (2025-12)
# Append a new contact-status row (row 5) to the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 5

$ws.Cells.Item($row, 1).Value = "DNET COMMUNICATIONS"

# BBM_STD (col B) has no value for this record - write it as an explicit
# empty text string (matching the other rows) via the quote-prefix trick,
# then reset the style so no extra formatting sticks.
$ws.Cells.Item($row, 2).Value = "'"
$ws.Cells.Item($row, 2).Style = "Normal"

$ws.Cells.Item($row, 3).Value = "OS"

# ACCOUNT_NO looks numeric but must stay text - force text storage with a
# leading apostrophe, then reset the style so no extra formatting sticks.
$ws.Cells.Item($row, 4).Value = "'9028244416"
$ws.Cells.Item($row, 4).Style = "Normal"

$ws.Cells.Item($row, 5).Value = "2025-12-02 14:28"

# LAST_WHATSAPP_TIME (col F) likewise has no value for this record.
$ws.Cells.Item($row, 6).Value = "'"
$ws.Cells.Item($row, 6).Style = "Normal"

$ws.Cells.Item($row, 7).Value = "2025-12"
